$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the KNeighborsRegressor block (original rows 51-58); everything below
# shifts up by 8 rows, matching the new dimension A1:F66.
$ws.Range("51:58").Delete()
